$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 189, pushing existing rows 189-287 down to 192-290.
$ws.Rows("189:191").Insert()

# Shared values for the three newly inserted rows (same market/product metadata
# as every other data row in this sheet).
$commonA = 3
$commonB = "Femacal de La Calera"
$commonC = "Coquimbo"
$commonE = 5
$commonF = "Fruta"
$commonG = 100107
$commonH = "Otros"
$commonI = 100107002
$commonJ = "Chirimoya"
$commonK = "Cultivar IV Región"
$commonQ = "$/bandeja 10 kilos"
$commonR = "Provincia del Elquí"
$commonT = 10
$newDate = 45126

function Set-DataRow($r, $quality, $volumen, $precio, $precioKg) {
    $ws.Cells.Item($r, 1).Value = $commonA
    $ws.Cells.Item($r, 2).Value = $commonB
    $ws.Cells.Item($r, 3).Value = $commonC
    $ws.Cells.Item($r, 4).Value = $newDate
    $ws.Cells.Item($r, 5).Value = $commonE
    $ws.Cells.Item($r, 6).Value = $commonF
    $ws.Cells.Item($r, 7).Value = $commonG
    $ws.Cells.Item($r, 8).Value = $commonH
    $ws.Cells.Item($r, 9).Value = $commonI
    $ws.Cells.Item($r, 10).Value = $commonJ
    $ws.Cells.Item($r, 11).Value = $commonK
    $ws.Cells.Item($r, 12).Value = $quality
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $precio
    $ws.Cells.Item($r, 15).Value = $precio
    $ws.Cells.Item($r, 16).Value = $precio
    $ws.Cells.Item($r, 17).Value = $commonQ
    $ws.Cells.Item($r, 18).Value = $commonR
    $ws.Cells.Item($r, 19).Value = $precioKg
    $ws.Cells.Item($r, 20).Value = $commonT
}

Set-DataRow 189 "Especial" 40 32000 3200
Set-DataRow 190 "Primera"  38 30000 3000
Set-DataRow 191 "Segunda"  35 27000 2700
